$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 3: a timedelta-style elapsed-time value in A3, plain number in B3 ---
$ws.Range("A3").NumberFormat = "[hh]:mm:ss"
$ws.Range("A3").Value2 = 10.6320601851852
$ws.Range("B3").Value2 = 17

# --- B1 / B2 pick up an explicit "General" number format (new style) ---
$ws.Range("B1:B2").NumberFormat = "General"

# --- Move the selection to A4, the cell below the new data ---
$ws.Range("A4").Select()

# --- Header/footer font label ("Regular" -> "Обычный") ---
$ws.PageSetup.CenterHeader = '&C&"Times New Roman,Обычный"&12&A'
$ws.PageSetup.CenterFooter = '&C&"Times New Roman,Обычный"&12Page &P'
